$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column for rows 2-230 from 45177 to 45178
$ws.Range("C2:C230").Value = 45178
